# Fix Package email batch detection: add delay after batch marker creation
# to prevent race condition. This re-ran the BOP smoke-test suite; the
# Summary sheet and the per-scenario timing sheets reflect the new run.

function Set-TextValue {
    # Assigns $Text to $Range while forcing a Text cell (so numeric-looking
    # strings like quote numbers / durations aren't coerced into numbers),
    # and restores the cell's original style afterward so no stray
    # "quote prefix" formatting is left behind.
    param($Range, [string]$Text)
    $origStyle = $Range.Style
    $Range.Value = "'" + $Text
    $Range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

# Row 2: BOP (OH) -> BOP (WI)
$summary.Range("B2").Value = "BOP (WI)"
Set-TextValue $summary.Range("F2") "222.81"
$summary.Range("G2").Value = "WI"

# Row 3: BOP (WI) -> BOP (DE)
$summary.Range("B3").Value = "BOP (DE)"
Set-TextValue $summary.Range("C3") "3003179622"
Set-TextValue $summary.Range("D3") "1003052773"
Set-TextValue $summary.Range("F3") "464.79"
$summary.Range("G3").Value = "DE"

# Row 4: BOP (PA) -> BOP (WI)
$summary.Range("B4").Value = "BOP (WI)"
Set-TextValue $summary.Range("C4") "3003179624"
Set-TextValue $summary.Range("D4") "1003052774"
Set-TextValue $summary.Range("F4") "486.61"
$summary.Range("G4").Value = "WI"

# Row 5: BOP (MI) now fails
$summary.Range("C5").Value = "N/A"
$summary.Range("D5").Value = "N/A"
$summary.Range("E5").Value = "FAILED"
Set-TextValue $summary.Range("F5") "107.42"

# Row 6: BOP (DE) -> BOP (OH)
$summary.Range("B6").Value = "BOP (OH)"
Set-TextValue $summary.Range("C6") "3003179627"
Set-TextValue $summary.Range("D6") "1003052775"
Set-TextValue $summary.Range("F6") "455.67"
$summary.Range("G6").Value = "OH"

# ---------------------------------------------------------------------
# BOP_1 sheet: last milestone row (the failure) disappears and the
# remaining two rows pick up the re-run's durations/timestamps.
# ---------------------------------------------------------------------
$bop1 = $wb.Worksheets.Item("BOP_1")
$bop1.Range("C2").Value = "68.82s"
$bop1.Range("D2").Value = "2025-12-24T18:53:11.574Z"
$bop1.Range("C3").Value = "153.99s"
$bop1.Range("D3").Value = "2025-12-24T18:55:45.570Z"
$bop1.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# BOP_2 sheet: all 6 milestone rows re-timed.
# ---------------------------------------------------------------------
$bop2 = $wb.Worksheets.Item("BOP_2")
$bop2.Range("C2").Value = "72.13s"
$bop2.Range("D2").Value = "2025-12-24T18:53:11.104Z"
$bop2.Range("C3").Value = "135.97s"
$bop2.Range("D3").Value = "2025-12-24T18:55:27.079Z"
$bop2.Range("C4").Value = "22.91s"
$bop2.Range("D4").Value = "2025-12-24T18:55:49.996Z"
$bop2.Range("C5").Value = "34.73s"
$bop2.Range("D5").Value = "2025-12-24T18:56:24.724Z"
$bop2.Range("C6").Value = "50.86s"
$bop2.Range("D6").Value = "2025-12-24T18:57:15.581Z"
$bop2.Range("C7").Value = "148.19s"
$bop2.Range("D7").Value = "2025-12-24T18:59:43.776Z"

# ---------------------------------------------------------------------
# BOP_3 sheet: all 6 milestone rows re-timed.
# ---------------------------------------------------------------------
$bop3 = $wb.Worksheets.Item("BOP_3")
$bop3.Range("C2").Value = "68.82s"
$bop3.Range("D2").Value = "2025-12-24T18:53:11.574Z"
$bop3.Range("C3").Value = "153.99s"
$bop3.Range("D3").Value = "2025-12-24T18:55:45.570Z"
$bop3.Range("C4").Value = "33.64s"
$bop3.Range("D4").Value = "2025-12-24T18:56:19.215Z"
$bop3.Range("C5").Value = "35.54s"
$bop3.Range("D5").Value = "2025-12-24T18:56:54.758Z"
$bop3.Range("C6").Value = "44.05s"
$bop3.Range("D6").Value = "2025-12-24T18:57:38.809Z"
$bop3.Range("C7").Value = "150.57s"
$bop3.Range("D7").Value = "2025-12-24T19:00:09.381Z"

# ---------------------------------------------------------------------
# BOP_4 sheet: run now fails at the first milestone; rows 3-7 disappear
# and row 2 becomes the failure record.
# ---------------------------------------------------------------------
$bop4 = $wb.Worksheets.Item("BOP_4")
$bop4.Range("A2").Value = "Test Execution Failed"
$bop4.Range("B2").Value = "FAILED"
$bop4.Range("C2").Value = "107.42s"
$bop4.Range("D2").Value = "2025-12-24T19:02:02.109Z"
$bop4.Rows("3:7").Delete()

# ---------------------------------------------------------------------
# BOP_5 sheet: all 6 milestone rows re-timed.
# ---------------------------------------------------------------------
$bop5 = $wb.Worksheets.Item("BOP_5")
$bop5.Range("C2").Value = "69.80s"
$bop5.Range("D2").Value = "2025-12-24T19:01:24.529Z"
$bop5.Range("C3").Value = "133.29s"
$bop5.Range("D3").Value = "2025-12-24T19:03:37.820Z"
$bop5.Range("C4").Value = "27.18s"
$bop5.Range("D4").Value = "2025-12-24T19:04:04.998Z"
$bop5.Range("C5").Value = "35.69s"
$bop5.Range("D5").Value = "2025-12-24T19:04:40.689Z"
$bop5.Range("C6").Value = "41.07s"
$bop5.Range("D6").Value = "2025-12-24T19:05:21.762Z"
$bop5.Range("C7").Value = "148.64s"
$bop5.Range("D7").Value = "2025-12-24T19:07:50.402Z"
